# Updates the cryptos list worksheet with refreshed prices / 1h volume
# percentages (and a rank swap between Hedera and PancakeSwap), matching
# the automated "Updated cryptos list ... with GitHub Actions" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a value while forcing it to be stored as literal text.
# Several "Price" values (column D) are plain decimal numbers (e.g.
# "0.487", "21.35"); Excel's COM layer would otherwise auto-convert
# those into numeric cells. Temporarily switching the cell to Text
# number format preserves the original string, and restoring the prior
# Style afterwards avoids leaving any stray formatting behind.
function Set-TextValue {
    param($range, [string]$value)
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

# Row 2 - Bitcoin
Set-TextValue $ws.Range("D2") "26.820.95"
$ws.Range("E2").Value = "  +0.09%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.541.88"
$ws.Range("E3").Value = "  -1.66%  "

# Row 5 - BNB
Set-TextValue $ws.Range("D5") "205.78"
$ws.Range("E5").Value = "  -0.43%  "

# Row 6 - XRP
Set-TextValue $ws.Range("D6") "0.487"
$ws.Range("E6").Value = "  -0.65%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.22%  "

# Row 8 - Cardano
$ws.Range("E8").Value = "  -0.38%  "

# Row 9 - Solana
Set-TextValue $ws.Range("D9") "21.35"
$ws.Range("E9").Value = "  -3.00%  "

# Row 10 - Dogecoin
Set-TextValue $ws.Range("D10") "0.0582"
$ws.Range("E10").Value = "  -0.54%  "

# Row 11 - TRON
Set-TextValue $ws.Range("D11") "0.0855"
$ws.Range("E11").Value = "  -1.10%  "

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "1.761.05"
$ws.Range("E12").Value = "  -1.71%  "

# Row 13 - WrappedEther
$ws.Range("D13").Value = "1.550.37"
$ws.Range("E13").Value = "  -1.07%  "

# Row 14 - Polkadot
$ws.Range("E14").Value = "  -1.55%  "

# Row 15 - Polygon
$ws.Range("E15").Value = "  -1.08%  "

# Row 16 - WrappedBTC
$ws.Range("D16").Value = "26.812.38"
$ws.Range("E16").Value = "  +0.05%  "

# Row 17 - Litecoin
Set-TextValue $ws.Range("D17") "61.22"
$ws.Range("E17").Value = "  -0.33%  "

# Row 18 - BitcoinCash
Set-TextValue $ws.Range("D18") "214.53"
$ws.Range("E18").Value = "  -0.26%  "

# Row 19 - Chainlink
$ws.Range("E19").Value = "  -2.51%  "

# Row 21 - Dai
$ws.Range("E21").Value = "  +0.19%  "

# Row 22 - Uniswap
$ws.Range("E22").Value = "  -3.14%  "

# Row 23 - Avalanche
Set-TextValue $ws.Range("D23") "9.16"
$ws.Range("E23").Value = "  -1.62%  "

# Row 24 - Toncoin
Set-TextValue $ws.Range("D24") "1.93"
$ws.Range("E24").Value = "  -3.19%  "

# Row 25 - Monero
Set-TextValue $ws.Range("D25") "152.30"
$ws.Range("E25").Value = "  -0.50%  "

# Row 26 - Cosmos
Set-TextValue $ws.Range("D26") "6.60"
$ws.Range("E26").Value = "  -2.22%  "

# Row 27 - EthereumClassic
Set-TextValue $ws.Range("D27") "14.82"
$ws.Range("E27").Value = "  -1.03%  "

# Row 28 - BinanceUSD
$ws.Range("E28").Value = "  +0.20%  "

# Row 29 - Stellar
$ws.Range("E29").Value = "  -0.90%  "

# Rows 30/31 - Hedera and PancakeSwap swap rank positions
$ws.Range("B30").Value = "Hedera"
$ws.Range("C30").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue $ws.Range("D30") "0.0458"
$ws.Range("E30").Value = "  -1.80%  "

$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue $ws.Range("D31") "1.10"
$ws.Range("E31").Value = "  -1.69%  "

# Row 32 - Filecoin
$ws.Range("E32").Value = "  +1.83%  "

# Row 33 - Maker
$ws.Range("D33").Value = "1.367.53"
$ws.Range("E33").Value = "  -2.10%  "

# Row 34 - InternetComputer(DFINITY)
$ws.Range("E34").Value = "  -0.02%  "

# Row 35 - LidoDAOToken
$ws.Range("E35").Value = "  -1.57%  "

# Row 36 - TrustWalletToken
Set-TextValue $ws.Range("D36") "0.965"

# Row 37 - HuobiToken
$ws.Range("E37").Value = "  +0.09%  "

# Row 38 - VeChain
$ws.Range("E38").Value = "  +0.73%  "

# Row 39 - ImmutableX
$ws.Range("E39").Value = "  -2.03%  "

# Row 40 - ARBITRUM
$ws.Range("E40").Value = "  -1.35%  "

# Row 41 - FraxShare
$ws.Range("E41").Value = "  +7.89%  "

# Row 42 - WEMIXToken
Set-TextValue $ws.Range("D42") "0.991"
$ws.Range("E42").Value = "  +0.41%  "

# Row 43 - MXToken
$ws.Range("E43").Value = "  +0.73%  "

# Row 44 - Aave
Set-TextValue $ws.Range("D44") "63.02"
$ws.Range("E44").Value = "  -0.55%  "

# Row 45 - RenderToken
Set-TextValue $ws.Range("D45") "1.74"
$ws.Range("E45").Value = "  -3.80%  "

# Row 46 - RocketPoolETH
$ws.Range("D46").Value = "1.675.53"
$ws.Range("E46").Value = "  -1.69%  "

# Row 48 - Cronos
$ws.Range("E48").Value = "  +3.47%  "

# Row 49 - BabyDogeCoin
$ws.Range("D49").Value = "0.0₇0974"
$ws.Range("E49").Value = "  -1.17%  "

# Row 51 - Algorand
$ws.Range("E51").Value = "  -1.49%  "
